$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.441.30'
$ws.Range('E2').Value = '  -1.08%  '

$ws.Range('D3').Value = '1.735.26'
$ws.Range('E3').Value = '  -1.24%  '

$ws.Range('E4').Value = '  -0.38%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '322.58'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4531'
$ws.Range('E7').Value = '  +6.65%  '

$ws.Range('E8').Value = '  -3.02%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07366'
$ws.Range('E9').Value = '  -2.30%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.30'
$ws.Range('E10').Value = '  -2.96%  '

$ws.Range('E11').Value = '  -1.77%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.51%  '

$ws.Range('E13').Value = '  -1.22%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.916'
$ws.Range('E14').Value = '  -2.03%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.070'
$ws.Range('E15').Value = '  -2.70%  '

$ws.Range('D16').Value = '1.729.12'
$ws.Range('E16').Value = '  -3.18%  '

$ws.Range('E17').Value = '  -0.43%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001051'
$ws.Range('E18').Value = '  -2.11%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06334'
$ws.Range('E19').Value = '  -0.82%  '

$ws.Range('E20').Value = '  -0.35%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.60'
$ws.Range('E21').Value = '  -2.63%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.738'
$ws.Range('E22').Value = '  -2.78%  '

$ws.Range('D23').Value = '27.477.16'
$ws.Range('E23').Value = '  -1.13%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.11'
$ws.Range('E24').Value = '  -0.68%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.062'
$ws.Range('E25').Value = '  -1.85%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '161.42'
$ws.Range('E26').Value = '  +0.55%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.90'
$ws.Range('E27').Value = '  -2.02%  '

$ws.Range('D28').Value = '1.926.39'
$ws.Range('E28').Value = '  -2.75%  '

$ws.Range('E29').Value = '  -4.23%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.57'
$ws.Range('E30').Value = '  -0.49%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.041'
$ws.Range('E31').Value = '  -6.63%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09118'
$ws.Range('E32').Value = '  +2.44%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.653'
$ws.Range('E33').Value = '  -0.27%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.381'
$ws.Range('E34').Value = '  -3.41%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02265'
$ws.Range('E35').Value = '  -1.48%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.57'
$ws.Range('E36').Value = '  -5.33%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05940'
$ws.Range('E37').Value = '  -1.54%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2043'
$ws.Range('E38').Value = '  -2.97%  '

$ws.Range('E39').Value = '  -1.59%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.855'
$ws.Range('E40').Value = '  -2.20%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.190'
$ws.Range('E41').Value = '  +0.37%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.368'
$ws.Range('E42').Value = '  -1.88%  '

$ws.Range('E43').Value = '  -2.52%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.07'
$ws.Range('E44').Value = '  -1.71%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.697'
$ws.Range('E45').Value = '  -0.07%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5785'
$ws.Range('E46').Value = '  -1.42%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '121.72'
$ws.Range('E47').Value = '  -1.29%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.922'
$ws.Range('E48').Value = '  -3.25%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06834'
$ws.Range('E49').Value = '  +0.13%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.109'
$ws.Range('E50').Value = '  -5.57%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '70.82'
$ws.Range('E51').Value = '  -3.73%  '
